$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the title
#    paragraph (paragraph 1). To avoid a save-time quirk where changing a
#    freshly-inserted paragraph's style away from a Heading style causes it
#    to vanish, we build the paragraph next to an existing Normal-styled
#    paragraph (paragraph 3) and then Cut/Paste it into its final position.
# ---------------------------------------------------------------------------

$normalAnchor = $d.Paragraphs.Item(3)
$anchorRange = $normalAnchor.Range
$anchorRange.Collapse(0)             # wdCollapseEnd
$anchorRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(4)
$newStart = $newPara.Range.Start
$insPoint = $d.Range($newStart, $newStart)
$insPoint.InsertAfter("Meta description: Try Arctic Valor slot for free and read our review. A high RTP and interesting gameplay make this slot worth playing.")

# Bold only the "Meta description" label (16 characters)
$boldRange = $d.Range($newStart, $newStart + 16)
$boldRange.Font.Bold = 1

# Cut the whole new paragraph (text + its trailing paragraph mark) ...
$newPara2 = $d.Paragraphs.Item(4)
$cutRange = $d.Range($newPara2.Range.Start, $newPara2.Range.End)
$cutRange.Cut()

# ... and paste it right after the title paragraph (paragraph 1).
$titlePara = $d.Paragraphs.Item(1)
$pastePoint = $d.Range($titlePara.Range.End, $titlePara.Range.End)
$pastePoint.Paste()

# ---------------------------------------------------------------------------
# 2. Remove the duplicate bold title paragraph that used to sit just before
#    the closing italic meta-description paragraph at the end of the doc.
# ---------------------------------------------------------------------------

$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($count - 1)
$dupRange = $d.Range($dupTitlePara.Range.Start, $dupTitlePara.Range.End)
$dupRange.Delete()

# ---------------------------------------------------------------------------
# 3. Replace the text of the final (italic) paragraph with the new prompt
#    text, keeping its existing italic formatting intact. The search is
#    scoped to just that final paragraph's range so the earlier, similarly
#    worded "Meta description" paragraph is left untouched.
# ---------------------------------------------------------------------------

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$searchRange = $lastPara.Range
$find = $searchRange.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
  "Try Arctic Valor slot for free and read our review. A high RTP and interesting gameplay make this slot worth playing.",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "Prompt: Create a feature image for Arctic Valor that features a happy Maya warrior with glasses in a cartoon style. The main colors used should be blue and white to match the icy theme of the game. The warrior should be holding a shield with a precious gemstone at the center, and in the background, there should be swirling snow and icicles hanging from the top. The image should be action-packed and showcase the excitement of the game.",
  2
)

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
